$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId=1 / index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 237
$ws1.Range("F4").Value = 214
$ws1.Range("F6").Value = 650
$ws1.Range("F8").Value = 438
$ws1.Range("F9").Value = 4227
$ws1.Range("F11").Value = 444
$ws1.Range("F13").Value = 972
$ws1.Range("F16").Value = 1892
$ws1.Range("F17").Value = 2937
$ws1.Range("F18").Value = 1769
$ws1.Range("F19").Value = 103
$ws1.Range("F21").Value = 160
$ws1.Range("F23").Value = 911
$ws1.Range("F24").Value = 284
$ws1.Range("F26").Value = 2225
$ws1.Range("F27").Value = 975
$ws1.Range("F28").Value = 2300
$ws1.Range("F30").Value = 679
$ws1.Range("F31").Value = 478
$ws1.Range("F33").Value = 872
$ws1.Range("F34").Value = 395
$ws1.Range("F35").Value = 1057
$ws1.Range("F36").Value = 876
$ws1.Range("F37").Value = 1141
$ws1.Range("F38").Value = 4
$ws1.Range("F39").Value = 320
$ws1.Range("F40").Value = 497
$ws1.Range("F41").Value = 352
$ws1.Range("F42").Value = 268
$ws1.Range("F43").Value = 3450

# Sheet "演出" (sheetId=2 / index 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 873
$ws2.Range("F17").Value = 1
$ws2.Range("F22").Value = 25
$ws2.Range("F23").Value = 1

# Sheet "全部类型" (sheetId=4 / index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 237
$ws4.Range("F4").Value = 214
$ws4.Range("F7").Value = 650
$ws4.Range("F9").Value = 438
$ws4.Range("F10").Value = 4228
$ws4.Range("F16").Value = 2937
$ws4.Range("F18").Value = 1769
$ws4.Range("F19").Value = 103
$ws4.Range("F22").Value = 160
$ws4.Range("F23").Value = 873
$ws4.Range("F27").Value = 911
$ws4.Range("F28").Value = 284
$ws4.Range("F29").Value = 2225
$ws4.Range("F32").Value = 975
$ws4.Range("F33").Value = 2300
$ws4.Range("F34").Value = 679
$ws4.Range("F35").Value = 478
$ws4.Range("F36").Value = 872
$ws4.Range("F37").Value = 395
$ws4.Range("F38").Value = 1057
$ws4.Range("F39").Value = 876
$ws4.Range("F40").Value = 1141
$ws4.Range("F41").Value = 320
$ws4.Range("F42").Value = 497
$ws4.Range("F44").Value = 352
$ws4.Range("F46").Value = 25
$ws4.Range("F47").Value = 268
$ws4.Range("F48").Value = 3450
